$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original row values for rows 3, 4 and 5 for the columns that change
$cols = @("A","B","D","E","F","G","H","P","Q","R","AI")

$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range("${col}3").Value()
    $row4[$col] = $ws.Range("${col}4").Value()
    $row5[$col] = $ws.Range("${col}5").Value()
}

# Cyclic rotation: row3 <- row4, row4 <- row5, row5 <- row3 (original)
foreach ($col in $cols) {
    $ws.Range("${col}3").Value = $row4[$col]
    $ws.Range("${col}4").Value = $row5[$col]
    $ws.Range("${col}5").Value = $row3[$col]
}
